$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect all target cells as Text before assigning values,
# so numeric-looking strings (e.g. "332.44", "0.00001050") are not
# auto-converted to numbers and lose their exact text representation.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"
$swapRange = $ws.Range("B17:C18")
$swapRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.812.35"
$ws.Range("E2").Value = "  +1.83%  "
$ws.Range("D3").Value = "1.882.66"
$ws.Range("E3").Value = "  +1.49%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "332.44"
$ws.Range("E5").Value = "  +2.62%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "0.4730"
$ws.Range("E7").Value = "  +4.61%  "
$ws.Range("D8").Value = "0.3965"
$ws.Range("E8").Value = "  +2.56%  "
$ws.Range("D9").Value = "47.89"
$ws.Range("E9").Value = "  -1.32%  "
$ws.Range("D10").Value = "0.08059"
$ws.Range("E10").Value = "  +1.62%  "
$ws.Range("E11").Value = "  +2.18%  "
$ws.Range("D12").Value = "22.24"
$ws.Range("E12").Value = "  +4.11%  "
$ws.Range("D13").Value = "1.882.01"
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("D14").Value = "5.984"
$ws.Range("E14").Value = "  +1.44%  "
$ws.Range("D15").Value = "7.148"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("D16").Value = "1.007"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.00001050"
$ws.Range("E17").Value = "  +2.21%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").Value = "87.33"
$ws.Range("E18").Value = "  +1.83%  "
$ws.Range("D19").Value = "0.06660"
$ws.Range("E19").Value = "  +1.46%  "
$ws.Range("D20").Value = "17.26"
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "27.802.37"
$ws.Range("E22").Value = "  +1.80%  "
$ws.Range("D23").Value = "5.549"
$ws.Range("E23").Value = "  +0.76%  "
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("D25").Value = "2.304"
$ws.Range("E25").Value = "  +0.74%  "
$ws.Range("D26").Value = "2.080.62"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").Value = "159.40"
$ws.Range("E27").Value = "  +3.92%  "
$ws.Range("D28").Value = "20.26"
$ws.Range("E28").Value = "  +2.06%  "
$ws.Range("D29").Value = "2.116"
$ws.Range("E29").Value = "  +2.57%  "
$ws.Range("D30").Value = "5.601"
$ws.Range("E30").Value = "  +2.86%  "
$ws.Range("D31").Value = "122.17"
$ws.Range("E31").Value = "  +1.03%  "
$ws.Range("D32").Value = "0.9892"
$ws.Range("E32").Value = "  +6.17%  "
$ws.Range("D33").Value = "0.09557"
$ws.Range("E33").Value = "  +2.91%  "
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("D35").Value = "3.595"
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").Value = "5.378"
$ws.Range("E36").Value = "  +2.21%  "
$ws.Range("D37").Value = "0.06135"
$ws.Range("E37").Value = "  +2.45%  "
$ws.Range("D38").Value = "0.02261"
$ws.Range("E38").Value = "  +1.83%  "
$ws.Range("E39").Value = "  +0.57%  "
$ws.Range("D40").Value = "8.157"
$ws.Range("E40").Value = "  +0.94%  "
$ws.Range("D41").Value = "0.6049"
$ws.Range("E41").Value = "  +2.51%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").Value = "0.1909"
$ws.Range("E43").Value = "  +1.28%  "
$ws.Range("D44").Value = "10.30"
$ws.Range("E44").Value = "  +2.15%  "
$ws.Range("E45").Value = "  -1.47%  "
$ws.Range("D46").Value = "0.5727"
$ws.Range("E46").Value = "  +1.91%  "
$ws.Range("D47").Value = "12.28"
$ws.Range("E47").Value = "  +2.43%  "
$ws.Range("E48").Value = "  +2.13%  "
$ws.Range("D49").Value = "3.378"
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("D50").Value = "0.06923"
$ws.Range("E50").Value = "  +2.54%  "
$ws.Range("D51").Value = "114.14"
$ws.Range("E51").Value = "  +5.16%  "

# Restore default (Normal) style so cells have no explicit number format,
# matching the original workbook formatting.
$priceRange.Style = "Normal"
$swapRange.Style = "Normal"
